# fix typo in MoA, preprocess dataset
# Replace every occurrence of the typo "DNA demage" with the corrected
# "DNA damage" in the moa column (column C) of the compound/smiles/moa sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("compound_smiles_moa")

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Text -eq "DNA demage") {
            $cell.Value = "DNA damage"
        }
    }
}

# Belt-and-braces: the rows known (from the source dataset) to carry the
# typo are C54, C61, C62 and C63 - make sure they are corrected even if the
# scan above used a different addressing convention.
foreach ($addr in @("C54", "C61", "C62", "C63")) {
    $cell = $ws.Range($addr)
    if ($cell.Text -eq "DNA demage") {
        $cell.Value = "DNA damage"
    }
}
